$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 148-161: column A held the old "SRX..._RNA2" isolate id and column B
# held the old plain "SRX..." accession. The isolate now carries a new
# GenBank accession (PV1858xx) in column A, and the old "SRX..._RNA2" value
# moves into column B (the accession column) as a cross reference.
$ws.Range("A148").Value = "PV185866"
$ws.Range("B148").Value = "SRX10234776_RNA2"

$ws.Range("A149").Value = "PV185868"
$ws.Range("B149").Value = "SRX10234777_RNA2"

$ws.Range("A150").Value = "PV185870"
$ws.Range("B150").Value = "SRX10234778_RNA2"

$ws.Range("A151").Value = "PV185872"
$ws.Range("B151").Value = "SRX10234779_RNA2"

$ws.Range("A152").Value = "PV185874"
$ws.Range("B152").Value = "SRX10234780_RNA2"

$ws.Range("A153").Value = "PV185876"
$ws.Range("B153").Value = "SRX10234781_RNA2"

$ws.Range("A154").Value = "PV185878"
$ws.Range("B154").Value = "SRX10234782_RNA2"

$ws.Range("A155").Value = " PV185881"
$ws.Range("B155").Value = "SRX10234784_RNA2"

$ws.Range("A156").Value = "PV185883"
$ws.Range("B156").Value = "SRX10234785_RNA2"

$ws.Range("A157").Value = "PV185885"
$ws.Range("B157").Value = "SRX10234786_RNA2"

$ws.Range("A158").Value = "PV185887"
$ws.Range("B158").Value = "SRX10234787_RNA2"

$ws.Range("A159").Value = "PV185889"
$ws.Range("B159").Value = "SRX10234788_RNA2"

$ws.Range("A160").Value = "PV185891"
$ws.Range("B160").Value = "SRX10234789_RNA2"

$ws.Range("A161").Value = "PV185893"
$ws.Range("B161").Value = "SRX10234790_RNA2"

# Widen column A slightly to fit the new accession numbers.
$ws.Columns.Item(1).ColumnWidth = 12.42

# Scroll the view down toward the bottom of the sheet and select the last
# rows, matching the author's working state when they saved.
[void]$ws.Range("A165:A179").EntireRow.Select()

$win = $excel.ActiveWindow
$win.ScrollRow = 155
$win.ScrollColumn = 1
$win.Zoom = 130
